# Update "limitele de incredere" (confidence limits) on the "Validare" sheet:
# add new Abatere / tp NN% columns (H:M) and refresh the f 1% / f 5% values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Validare")

# The sheet has no cell styles yet; touching .Style establishes a default
# style entry so subsequent value writes succeed.
$ws.Range("H1:M2").Style = "Normal"

# Row 2 values are stored as text (not numbers) in this workbook, so force a
# text number format on the numeric-looking cells before writing them.
$ws.Range("F2:M2").NumberFormat = "@"

# New header cells for the added confidence-limit columns
$ws.Range("H1").Value = "Abatere"
$ws.Range("I1").Value = "tp 80%"
$ws.Range("J1").Value = "tp 90%"
$ws.Range("K1").Value = "tp 95%"
$ws.Range("L1").Value = "tp 99%"
$ws.Range("M1").Value = "tp 99.9%"

# Updated / new values on row 2
$ws.Range("F2").Value = "13.27"
$ws.Range("G2").Value = "5.79"
$ws.Range("H2").Value = "304.613058889895"
$ws.Range("I2").Value = "1.44"
$ws.Range("J2").Value = "1.943"
$ws.Range("K2").Value = "2.447"
$ws.Range("L2").Value = "3.707"
$ws.Range("M2").Value = "5.959"
